# AF-611: tests for named ranges and area intersections are moved to
# temporary_excel_files. The B2 formula no longer references the undefined
# name "text" (which produced a #NAME? error); it now concatenates A1 and A2
# (both empty), which evaluates to an empty string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Formula = "=CONCATENATE(A1,A2)"
